# Refresh the cryptos price list (commit: "Updated cryptos list ... with GitHub
# Actions"). For every changed row we overwrite the Price (D) and Volume(1h) (E)
# cells with the new scraped values; rows 25/26 additionally had their Coin (B)
# and Link (C) cells swapped (Toncoin now ranks above LidoDAOToken).
#
# NOTE: column D holds prices as plain TEXT in this workbook (inline strings),
# not numbers -- several new values ("51.60", "1.000", "0.3778", ...) are
# numeric-looking, so Excel would otherwise auto-convert them into real numbers
# on assignment. Prefixing those with a leading apostrophe forces them to stay
# text, matching the original cell type/content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = '23.193.10'
$ws.Range("E2").Value = '  +0.29%  '

# Row 3 (Ethereum)
$ws.Range("D3").Value = '1.601.78'
$ws.Range("E3").Value = '  -0.05%  '

# Row 4 (TetherUSD)
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.07%  '

# Row 6 (BNB)
$ws.Range("D6").Value = '''303.08'
$ws.Range("E6").Value = '  +0.63%  '

# Row 7 (XRP)
$ws.Range("D7").Value = '''0.3778'
$ws.Range("E7").Value = '  -0.27%  '

# Row 8 (OKB)
$ws.Range("D8").Value = '''51.60'
$ws.Range("E8").Value = '  +3.12%  '

# Row 9 (Cardano)
$ws.Range("D9").Value = '''0.3614'
$ws.Range("E9").Value = '  -0.99%  '

# Row 10 (Polygon)
$ws.Range("D10").Value = '''1.266'
$ws.Range("E10").Value = '  +0.70%  '

# Row 11 (BinanceUSD)
$ws.Range("D11").Value = '''1.000'
$ws.Range("E11").Value = '  -0.13%  '

# Row 12 (Dogecoin)
$ws.Range("E12").Value = '  -0.09%  '

# Row 13 (Solana)
$ws.Range("D13").Value = '''22.59'
$ws.Range("E13").Value = '  -1.92%  '

# Row 14 (Polkadot)
$ws.Range("D14").Value = '''6.598'
$ws.Range("E14").Value = '  -0.05%  '

# Row 15 (Chainlink)
$ws.Range("D15").Value = '''7.392'
$ws.Range("E15").Value = '  -0.37%  '

# Row 16 (ShibaInu)
$ws.Range("D16").Value = '''0.00001248'
$ws.Range("E16").Value = '  -0.46%  '

# Row 17 (WrappedEther)
$ws.Range("D17").Value = '1.602.20'
$ws.Range("E17").Value = '  -0.10%  '

# Row 18 (Litecoin)
$ws.Range("D18").Value = '''93.66'
$ws.Range("E18").Value = '  +2.25%  '

# Row 19 (TRON)
$ws.Range("D19").Value = '''0.06869'
$ws.Range("E19").Value = '  +0.16%  '

# Row 20 (Avalanche)
$ws.Range("D20").Value = '''18.03'
$ws.Range("E20").Value = '  -1.36%  '

# Row 21 (Uniswap)
$ws.Range("D21").Value = '''6.532'
$ws.Range("E21").Value = '  -0.43%  '

# Row 22 (Dai)
$ws.Range("D22").Value = '''0.9994'
$ws.Range("E22").Value = '  -0.12%  '

# Row 23 (Cosmos)
$ws.Range("D23").Value = '''12.97'
$ws.Range("E23").Value = '  -0.36%  '

# Row 24 (WrappedBTC)
$ws.Range("D24").Value = '23.199.97'
$ws.Range("E24").Value = '  +0.29%  '

# Row 25 (Toncoin)
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '''2.393'
$ws.Range("E25").Value = '  +2.23%  '

# Row 26 (LidoDAOToken)
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '''2.997'
$ws.Range("E26").Value = '  +10.11%  '

# Row 27 (EthereumClassic)
$ws.Range("D27").Value = '''21.18'
$ws.Range("E27").Value = '  +0.42%  '

# Row 28 (Monero)
$ws.Range("D28").Value = '''150.06'
$ws.Range("E28").Value = '  -0.12%  '

# Row 29 (HuobiToken)
$ws.Range("D29").Value = '''5.233'
$ws.Range("E29").Value = '  -0.94%  '

# Row 30 (BitcoinCash)
$ws.Range("D30").Value = '''133.90'
$ws.Range("E30").Value = '  +1.52%  '

# Row 31 (WEMIXTOKEN)
$ws.Range("D31").Value = '''2.419'
$ws.Range("E31").Value = '  -0.21%  '

# Row 32 (Filecoin)
$ws.Range("D32").Value = '''6.817'
$ws.Range("E32").Value = '  -0.41%  '

# Row 33 (WrappedliquidstakedEther2.0)
$ws.Range("D33").Value = '1.780.34'
$ws.Range("E33").Value = '  +0.02%  '

# Row 34 (ImmutableX)
$ws.Range("D34").Value = '''0.9804'
$ws.Range("E34").Value = '  +3.66%  '

# Row 35 (Hedera)
$ws.Range("D35").Value = '''0.07554'
$ws.Range("E35").Value = '  -1.73%  '

# Row 36 (FraxShare)
$ws.Range("D36").Value = '''10.31'
$ws.Range("E36").Value = '  +2.92%  '

# Row 37 (VeChain)
$ws.Range("D37").Value = '''0.02725'
$ws.Range("E37").Value = '  -0.73%  '

# Row 38 (InternetComputer(DFINITY))
$ws.Range("D38").Value = '''6.138'
$ws.Range("E38").Value = '  -1.82%  '

# Row 39 (Algorand)
$ws.Range("D39").Value = '''0.2502'
$ws.Range("E39").Value = '  -1.64%  '

# Row 40 (Stellar)
$ws.Range("D40").Value = '''0.08794'
$ws.Range("E40").Value = '  -1.62%  '

# Row 41 (TheSandbox)
$ws.Range("D41").Value = '''0.7108'
$ws.Range("E41").Value = '  -0.02%  '

# Row 42 (TrustWalletToken)
$ws.Range("D42").Value = '''1.360'
$ws.Range("E42").Value = '  -2.07%  '

# Row 43 (Aptos)
$ws.Range("D43").Value = '''12.43'
$ws.Range("E43").Value = '  -2.26%  '

# Row 44 (EnergySwap)
$ws.Range("D44").Value = '''15.47'
$ws.Range("E44").Value = '  +0.14%  '

# Row 45 (Decentraland)
$ws.Range("D45").Value = '''0.6552'
$ws.Range("E45").Value = '  -1.05%  '

# Row 46 (NEARProtocol)
$ws.Range("D46").Value = '''2.309'
$ws.Range("E46").Value = '  +0.36%  '

# Row 47 (PancakeSwap)
$ws.Range("D47").Value = '''4.014'
$ws.Range("E47").Value = '  +1.00%  '

# Row 48 (Quant)
$ws.Range("D48").Value = '''132.32'
$ws.Range("E48").Value = '  +0.17%  '

# Row 49 (Cronos)
$ws.Range("D49").Value = '''0.07960'
$ws.Range("E49").Value = '  +0.31%  '

# Row 50 (Flow)
$ws.Range("D50").Value = '''1.207'
$ws.Range("E50").Value = '  -0.62%  '

# Row 51 (ThetaToken)
$ws.Range("D51").Value = '''1.230'
$ws.Range("E51").Value = '  +3.84%  '
